$d = $word.ActiveDocument

# --- Change 1: shorten the "to correct the problem of the atmosphere..." sentence ---
$d.Content.Find.Execute(
    "to correct the problem of the atmosphere. Basically, a sensor measures the distortion of the atmosphere every few milliseconds, and a mirror in the telescope deforms in order to adjust for those distortions.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "to correct the problem of the atmosphere using sensors and a deformable mirror.",
    2) | Out-Null

# --- Change 2: append the "Slide" outline paragraphs, and relocate the _GoBack bookmark ---

# Append an empty paragraph after the final "...are there any questions." paragraph.
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertParagraphAfter() | Out-Null

# Text for each new "Slide" paragraph.
$slides = @(
    "Slide 1: Introduce myself, steward observatory, mentor, project, what presentation will be about",
    "Slide 2: Adaptive optics introduction, what it is, why it’s necessary, what that means to my project",
    "Slide 3: Explain where the pictures were taken and what I was doing in the other picture. Dewar filling with liquid nitrogen and that the detector needed to be at a low temperature because of the nature of infrared photo taking",
    "Slide 4: Coding was focus, how calibrating the initial dark field set could lead to applied coefficients to other data sets to calibrate them as well",
    "Slide 5: What are ints? What are counts? What is linearity, and how it relates to ints and counts? What did I use to code? Which packages did I use?",
    "Slide 6: Process of code: read in the code, test orders to equations to fit, error plots",
    "Slide 7: Calibrated the counts to fourth order, for examples, look at these graphics!",
    "Slide 8: Explain pictures",
    "Slide 9: Impact on modern astronomy: EXOPLANETS BIIITCH"
)

foreach ($slide in $slides) {
    $r = $d.Paragraphs.Last.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter() | Out-Null
    $r2 = $d.Paragraphs.Last.Range
    $r2.InsertAfter($slide) | Out-Null
}

# Relocate the _GoBack bookmark: it used to sit between "correct" and "ed" in the
# "correction" paragraph; it now sits inside the new "Slide 5" paragraph, right
# after the word "What" that begins the second sentence ("What are counts?...").
$rng = $d.Content.Duplicate
$rng.Find.Execute("Slide 5: What", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rng) | Out-Null

Write-Output "done"
